# Applies the data refresh described in the commit "Tings are happening hard".
# Updates numeric length counters (M2_Len / FX_Len in columns C/D) and the
# date-serial columns (M2_1stDate/M2_LastDate/FX_1stDate/FX_LastDate in
# columns E/F/G/H) on the active worksheet to their refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 342
$ws.Range("F2").Value = 45444
$ws.Range("G2").Value = 30319
$ws.Range("H2").Value = 45505
$ws.Range("E3").Value = 30256
$ws.Range("F3").Value = 45444
$ws.Range("E4").Value = 30256
$ws.Range("F4").Value = 45444
$ws.Range("G4").Value = 30319
$ws.Range("H4").Value = 45505
$ws.Range("E5").Value = 30256
$ws.Range("F5").Value = 45444
$ws.Range("G5").Value = 30319
$ws.Range("H5").Value = 45505
$ws.Range("C6").Value = 451
$ws.Range("F6").Value = 45444
$ws.Range("G6").Value = 30319
$ws.Range("H6").Value = 45505
$ws.Range("E7").Value = 30256
$ws.Range("F7").Value = 45444
$ws.Range("G7").Value = 30319
$ws.Range("H7").Value = 45505
$ws.Range("D8").Value = 418
$ws.Range("E8").Value = 30256
$ws.Range("F8").Value = 45444
$ws.Range("H8").Value = 45505
$ws.Range("E9").Value = 30256
$ws.Range("F9").Value = 45444
$ws.Range("G9").Value = 30319
$ws.Range("H9").Value = 45505
$ws.Range("D10").Value = 491
$ws.Range("E10").Value = 30256
$ws.Range("F10").Value = 45444
$ws.Range("H10").Value = 45505
$ws.Range("E11").Value = 30225
$ws.Range("F11").Value = 45413
$ws.Range("G11").Value = 30319
$ws.Range("H11").Value = 45505
$ws.Range("C12").Value = 379
$ws.Range("D12").Value = 360
$ws.Range("F12").Value = 45444
$ws.Range("H12").Value = 45505
$ws.Range("C13").Value = 474
$ws.Range("F13").Value = 45413
$ws.Range("G13").Value = 30319
$ws.Range("H13").Value = 45505
$ws.Range("D14").Value = 404
$ws.Range("H14").Value = 45505
$ws.Range("G15").Value = 30286
$ws.Range("H15").Value = 45505
$ws.Range("C16").Value = 462
$ws.Range("D16").Value = 418
$ws.Range("F16").Value = 45413
$ws.Range("H16").Value = 45505
$ws.Range("C17").Value = 378
$ws.Range("D17").Value = 402
$ws.Range("F17").Value = 45444
$ws.Range("H17").Value = 45505
$ws.Range("D18").Value = 268
$ws.Range("E18").Value = 30256
$ws.Range("F18").Value = 45444
$ws.Range("H18").Value = 45505
$ws.Range("D19").Value = 406
$ws.Range("E19").Value = 30256
$ws.Range("F19").Value = 45444
$ws.Range("H19").Value = 45505
$ws.Range("C20").Value = 489
$ws.Range("F20").Value = 45444
$ws.Range("G20").Value = 30319
$ws.Range("H20").Value = 45505
$ws.Range("C21").Value = 317
$ws.Range("F21").Value = 45413
$ws.Range("G21").Value = 30319
$ws.Range("H21").Value = 45505
$ws.Range("C22").Value = 331
$ws.Range("D22").Value = 375
$ws.Range("F22").Value = 45444
$ws.Range("H22").Value = 45505
$ws.Range("C23").Value = 267
$ws.Range("D23").Value = 410
$ws.Range("F23").Value = 45383
$ws.Range("H23").Value = 45505
$ws.Range("C24").Value = 224
$ws.Range("D24").Value = 418
$ws.Range("F24").Value = 45474
$ws.Range("H24").Value = 45505
$ws.Range("D25").Value = 323
$ws.Range("E25").Value = 30225
$ws.Range("F25").Value = 45413
$ws.Range("H25").Value = 45505
$ws.Range("C26").Value = 341
$ws.Range("D26").Value = 321
$ws.Range("F26").Value = 45444
$ws.Range("H26").Value = 45505
$ws.Range("E27").Value = 30225
$ws.Range("F27").Value = 45413
$ws.Range("G27").Value = 30319
$ws.Range("H27").Value = 45505
$ws.Range("D28").Value = 388
$ws.Range("E28").Value = 30256
$ws.Range("F28").Value = 45444
$ws.Range("H28").Value = 45505
$ws.Range("C29").Value = 270
$ws.Range("D29").Value = 244
$ws.Range("F29").Value = 45444
$ws.Range("H29").Value = 45505
$ws.Range("D30").Value = 225
$ws.Range("E30").Value = 30225
$ws.Range("F30").Value = 45413
$ws.Range("H30").Value = 45505
$ws.Range("C31").Value = 402
$ws.Range("F31").Value = 45444
$ws.Range("G31").Value = 30319
$ws.Range("H31").Value = 45505
$ws.Range("G32").Value = 30319
$ws.Range("H32").Value = 45505
$ws.Range("C33").Value = 462
$ws.Range("D33").Value = 406
$ws.Range("F33").Value = 45444
$ws.Range("H33").Value = 45505
$ws.Range("C34").Value = 209
$ws.Range("D34").Value = 325
$ws.Range("F34").Value = 45413
$ws.Range("H34").Value = 45505
$ws.Range("C35").Value = 412
$ws.Range("D35").Value = 325
$ws.Range("F35").Value = 45413
$ws.Range("H35").Value = 45505
$ws.Range("D36").Value = 418
$ws.Range("E36").Value = 30256
$ws.Range("F36").Value = 45444
$ws.Range("H36").Value = 45505
$ws.Range("C37").Value = 474
$ws.Range("D37").Value = 325
$ws.Range("F37").Value = 45444
$ws.Range("H37").Value = 45505
$ws.Range("C38").Value = 367
$ws.Range("D38").Value = 375
$ws.Range("F38").Value = 45444
$ws.Range("H38").Value = 45505
$ws.Range("C39").Value = 235
$ws.Range("D39").Value = 231
$ws.Range("F39").Value = 45444
$ws.Range("H39").Value = 45505
$ws.Range("C40").Value = 294
$ws.Range("D40").Value = 322
$ws.Range("F40").Value = 45444
$ws.Range("H40").Value = 45505
$ws.Range("C41").Value = 402
$ws.Range("D41").Value = 244
$ws.Range("F41").Value = 45444
$ws.Range("H41").Value = 45505
$ws.Range("D42").Value = 226
$ws.Range("H42").Value = 45505
$ws.Range("C43").Value = 481
$ws.Range("D43").Value = 325
$ws.Range("F43").Value = 45444
$ws.Range("H43").Value = 45505
$ws.Range("C44").Value = 413
$ws.Range("D44").Value = 312
$ws.Range("F44").Value = 45413
$ws.Range("H44").Value = 45505
$ws.Range("D45").Value = 325
$ws.Range("H45").Value = 45505
$ws.Range("C46").Value = 342
$ws.Range("D46").Value = 306
$ws.Range("F46").Value = 45444
$ws.Range("H46").Value = 45505
$ws.Range("C47").Value = 343
$ws.Range("D47").Value = 268
$ws.Range("F47").Value = 45444
$ws.Range("H47").Value = 45505
$ws.Range("C48").Value = 365
$ws.Range("D48").Value = 323
$ws.Range("F48").Value = 45413
$ws.Range("H48").Value = 45505
$ws.Range("C49").Value = 305
$ws.Range("D49").Value = 321
$ws.Range("F49").Value = 45413
$ws.Range("H49").Value = 45505
$ws.Range("C50").Value = 366
$ws.Range("D50").Value = 246
$ws.Range("F50").Value = 45444
$ws.Range("H50").Value = 45505
$ws.Range("D51").Value = 325
$ws.Range("H51").Value = 45505
$ws.Range("D52").Value = 323
$ws.Range("H52").Value = 45505
